$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Offer Letter For: *" -> "Effective Date: *"
$ws.Range("C1").Value = "Effective Date: *"

# "Format Id:*" column is no longer used - clear it out entirely
$ws.Range("W1").ClearContents()

# Widen column H (Special Allowance) a bit
$ws.Columns.Item(8).ColumnWidth = 14.86

# Move the view back to the top-left and park the selection on F19
$ws.Range("A1").Select()
$ws.Range("F19").Select()
